$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3: new year headers 2020, 2021, 2022 (copy formatting from J3) ---
$ws.Range("J3").Copy()
$ws.Range("K3:M3").PasteSpecial(-4122)
$ws.Range("K3").Value = 2020
$ws.Range("L3").Value = 2021
$ws.Range("M3").Value = 2022

# --- Row 4: new data values (copy formatting from J4) ---
$ws.Range("J4").Copy()
$ws.Range("K4:M4").PasteSpecial(-4122)
$ws.Range("K4").Value = 308
$ws.Range("L4").Value = 212.1
$ws.Range("M4").Value = 723.8

# --- Row 5: fill E5:L5 with "-" using a new right-aligned style, M5 stays blank ---
$ws.Range("D5").Copy()
$ws.Range("E5:M5").PasteSpecial(-4122)
$ws.Range("E5:M5").HorizontalAlignment = -4152
$ws.Range("E5:L5").Value = "-"

# --- Row 6: new data values (copy formatting from I6), M6 stays blank ---
$ws.Range("I6").Copy()
$ws.Range("J6:M6").PasteSpecial(-4122)
$ws.Range("J6").Value = 9.8
$ws.Range("K6").Value = 9.8
$ws.Range("L6").Value = 9.8

# --- Row 7: new data values (copy formatting from I7) ---
$ws.Range("I7").Copy()
$ws.Range("J7:M7").PasteSpecial(-4122)
$ws.Range("J7").Value = 64
$ws.Range("K7").Value = 64
$ws.Range("L7").Value = 64
$ws.Range("M7").Value = 64

# --- Update selection to match target ---
$ws.Range("M14").Select()
